# Apply the odds updates described by the diff to Sheet1 (active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new numeric value
$changes = @{
    "F2"  = 1.6
    "L2"  = 1.39
    "P2"  = 1.87
    "Q2"  = 1.89
    "S2"  = 3.05
    "U2"  = 1.89

    "F3"  = 1.23
    "H3"  = 11.5
    "I3"  = 14
    "J3"  = 7.4
    "K3"  = 8.800000000000001
    "L3"  = 1.16
    "N3"  = 7.6
    "P3"  = 3.65
    "R3"  = 2.06
    "S3"  = 1.76
    "T3"  = 1.74
    "U3"  = 2.08
    "Y3"  = 970
    "AC3" = 21
    "AF3" = 14.5
    "AH3" = 30
    "AJ3" = 13.5
    "AK3" = 16.5
    "AL3" = 32
    "AM3" = 130
    "AN3" = 3.1
    "AO3" = 160

    "F4"  = 3.85
    "I4"  = 2.28
    "J4"  = 3.3
    "K4"  = 3.35
    "AF4" = 27
    "AH4" = 22
    "AO4" = 26

    "H5"  = 2.94
    "O5"  = 1.44
    "P5"  = 1.71
    "Q5"  = 2.32
    "R5"  = 1.27
    "Z5"  = 19
    "AA5" = 65
    "AE5" = 40
    "AI5" = 70
    "AK5" = 36
    "AL5" = 65

    "H6"  = 2.58
    "I6"  = 2.6
    "P6"  = 1.82
    "S6"  = 4
    "AG6" = 14
    "AH6" = 20
    "AI6" = 50
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
